$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column at D: this shifts old D (empty) -> E, old E (Precio por 2012 data)
# -> F, old F (Precio por 2016 data) -> G, old I/J (Media labels/formulas) -> J/K, etc.
# Excel auto-adjusts formula references (e.g. E33/B33 -> F33/B33) during the shift.
$ws.Columns("D").Insert()

# var 2012 / var 2016 labels (entered first so shared-string order matches
# the authored workbook: var 2012, var 2016, then the D/E headers, then sd)
$ws.Range("J6").Value = "var 2012"
$ws.Range("J7").Value = "var 2016"

# New column D header: "2012 por (precio - media_2012)^2"
$ws.Range("D1").Value = "2012 por (precio - media_2012)^2"
# New column E header: "2016 por (precio - media_2016)^2"
$ws.Range("E1").Value = "2016 por (precio - media_2016)^2"

# sd 2012 / sd 2016 labels
$ws.Range("J8").Value = "sd 2012"
$ws.Range("J9").Value = "sd 2016  "

# Fill D (2012 contribution to variance) and E (2016 contribution to
# variance) columns, row 2 set apart from 3:32 so the formula grouping
# mirrors the rest of the sheet (row 2 is always its own formula, 3:32
# share one formula group), matching how B*A / C*A were laid out.
$ws.Range("D2").Formula = "=B2*(A2-`$K`$4)^2"
$ws.Range("D3:D32").Formula = "=B3*(A3-`$K`$4)^2"

$ws.Range("E2").Formula = "=C2*(A2-`$K`$5)^2"
$ws.Range("E3:E32").Formula = "=C3*(A3-`$K`$5)^2"

# Re-assert G3:G32 (A*C, shifted from the original F column) as one range
# formula so it keeps a shared-formula grouping like the rest of the sheet
# (column insert alone leaves it as per-cell formulas).
$ws.Range("G3:G32").Formula = "=A3*C3"

# Row 33 totals for the new columns
$ws.Range("D33").Formula = "=SUM(D2:D32)"
$ws.Range("E33").Formula = "=SUM(E2:E32)"

# var 2012 / var 2016 formulas
$ws.Range("K6").Formula = "=D33/B33"
$ws.Range("K7").Formula = "=E33/C33"

# sd 2012 / sd 2016 formulas
$ws.Range("K8").Formula = "=SQRT(K6)"
$ws.Range("K9").Formula = "=SQRT(K7)"

# Selection/active cell ends on K8, matching the authored workbook state
$ws.Range("K8").Select()
